$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2-176
# from serial 45182 (2023-09-13) to serial 45184 (2023-09-15).
$ws.Range("C2:C176").Value = 45184
